# Add two new columns of room-seating documentation ("row_dist", "seat_dist")
# after the existing "room" column, shifting the old "seat_last" column from
# E to G.  Cell-entry order below matches the shared-string insertion order
# of the target workbook: header row first, then row 3 ("1,1"/"0,6"), then
# row 2 ("1,5" reused for both E2 and F2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("E1").Value = "row_dist"
$ws.Range("F1").Value = "seat_dist"
$ws.Range("G1").Value = "seat_last"

# Row 3
$ws.Range("E3").Value = "1,1"
$ws.Range("F3").Value = "0,6"
$ws.Range("G3").Value = "r2s3"

# Row 2
$ws.Range("E2").Value = "1,5"
$ws.Range("F2").Value = "1,5"
$ws.Range("G2").Value = "r2s7"

# Auto-size the two newly filled columns, as Excel does when data is typed
# into a previously-empty column.
$ws.Range("E:F").EntireColumn.AutoFit()

# Leave the selection where the user finished typing.
$ws.Range("E3:F3").Select()
